$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 83.58
$ws.Range("I15").Value = 83.58
$ws.Range("K15").Value = 250.74
$ws.Range("M15").Value = -81.74000000000001
$ws.Range("H19").Value = 1212280.2
$ws.Range("I19").Value = 1481587
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 1481587
$ws.Range("L19").Value = 400
$ws.Range("M19").Value = -1481412
$ws.Range("N19").Value = -750
$ws.Range("H61").Value = 270.9
$ws.Range("I61").Value = 245.44444
$ws.Range("K61").Value = 736.33332
$ws.Range("M61").Value = -564.33332
$ws.Range("H69").Value = 5128.25
$ws.Range("I69").Value = 5128.25
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 15384.75
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -14510.75
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 5128.25
$ws.Range("I72").Value = 5128.25
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 46154.25
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -41786.25
$ws.Range("N72").ClearContents()
$ws.Range("H112").Value = 10871178
$ws.Range("J112").Value = 1638.3334
$ws.Range("L112").Value = 4915.0002
$ws.Range("N112").Value = -7131.0002
$ws.Range("H124").Value = 42780
$ws.Range("J124").Value = 42780
$ws.Range("L124").Value = 42780
$ws.Range("N124").Value = -52600
$ws.Range("H129").Value = 898.5454999999999
$ws.Range("J129").Value = 1009.3333
$ws.Range("L129").Value = 3027.9999
$ws.Range("N129").Value = -13027.9999
$ws.Range("H137").Value = 2840.1785
$ws.Range("I137").Value = 2429.7173
$ws.Range("J137").Value = 4728.3
$ws.Range("K137").Value = 7289.151899999999
$ws.Range("L137").Value = 14184.9
$ws.Range("M137").Value = -4739.151899999999
$ws.Range("N137").Value = -19284.9
$ws.Range("H139").Value = 42510.477
$ws.Range("J139").Value = 42510.477
$ws.Range("L139").Value = 42510.477
$ws.Range("N139").Value = -52790.477

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7947.829
$ws.Range("I32").Value = 5894.691
$ws.Range("K32").Value = 5894.691
$ws.Range("M32").Value = -5607.691
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H122").Value = 2062.0356
$ws.Range("I122").Value = 1053.4762
$ws.Range("J122").Value = 5087.7144
$ws.Range("K122").Value = 3160.4286
$ws.Range("L122").Value = 15263.1432
$ws.Range("M122").Value = -710.4286000000002
$ws.Range("N122").Value = -20163.1432

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2038.3334
$ws.Range("I134").Value = 1182.7681
$ws.Range("J134").Value = 6957.8335
$ws.Range("K134").Value = 3548.3043
$ws.Range("L134").Value = 20873.5005
$ws.Range("M134").Value = -1013.3043
$ws.Range("N134").Value = -25943.5005
$ws.Range("H138").Value = 41507.5
$ws.Range("J138").Value = 41507.5
$ws.Range("L138").Value = 41507.5
$ws.Range("N138").Value = -51787.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2072.6616
$ws.Range("I31").Value = 856.9729599999999
$ws.Range("K31").Value = 856.9729599999999
$ws.Range("M31").Value = -561.9729599999999
$ws.Range("H34").Value = 2072.6616
$ws.Range("I34").Value = 856.9729599999999
$ws.Range("K34").Value = 856.9729599999999
$ws.Range("M34").Value = -654.9729599999999
$ws.Range("H58").Value = 1926.2667
$ws.Range("J58").Value = 5440
$ws.Range("L58").Value = 5440
$ws.Range("N58").Value = -5846
$ws.Range("H136").Value = 1926.2667
$ws.Range("J136").Value = 5440
$ws.Range("L136").Value = 16320
$ws.Range("N136").Value = -21420
$ws.Range("H138").Value = 39350
$ws.Range("J138").Value = 39350
$ws.Range("L138").Value = 39350
$ws.Range("N138").Value = -49630

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4827101.5
$ws.Range("J4").Value = 8661.333000000001
$ws.Range("L4").Value = 25983.999
$ws.Range("N4").Value = -26207.999
$ws.Range("H7").Value = 283.33334
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 375
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 1125
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -1349
$ws.Range("H14").Value = 123.318184
$ws.Range("I14").Value = 123.318184
$ws.Range("K14").Value = 369.954552
$ws.Range("M14").Value = -196.954552
$ws.Range("H137").Value = 4508.857
$ws.Range("J137").Value = 7516.5
$ws.Range("L137").Value = 22549.5
$ws.Range("N137").Value = -32749.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2849.0303
$ws.Range("I102").Value = 2280
$ws.Range("J102").Value = 3621.2856
$ws.Range("K102").Value = 2280
$ws.Range("L102").Value = 3621.2856
$ws.Range("M102").Value = -658
$ws.Range("N102").Value = -6865.2856
$ws.Range("H132").Value = 2511.5588
$ws.Range("I132").Value = 1531.28
$ws.Range("K132").Value = 4593.84
$ws.Range("M132").Value = -2063.84
$ws.Range("H138").Value = 52943.332
$ws.Range("J138").Value = 52943.332
$ws.Range("L138").Value = 52943.332
$ws.Range("N138").Value = -63223.332

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 835.4516
$ws.Range("I68").Value = 763.3
$ws.Range("K68").Value = 763.3
$ws.Range("M68").Value = -14.29999999999995
$ws.Range("H71").Value = 835.4516
$ws.Range("I71").Value = 763.3
$ws.Range("K71").Value = 3816.5
$ws.Range("M71").Value = -72.5
$ws.Range("H82").Value = 1058.9048
$ws.Range("I82").Value = 802.4666999999999
$ws.Range("J82").Value = 1700
$ws.Range("K82").Value = 802.4666999999999
$ws.Range("L82").Value = 1700
$ws.Range("M82").Value = -441.4666999999999
$ws.Range("N82").Value = -2422
$ws.Range("H85").Value = 1058.9048
$ws.Range("I85").Value = 802.4666999999999
$ws.Range("J85").Value = 1700
$ws.Range("K85").Value = 802.4666999999999
$ws.Range("L85").Value = 1700
$ws.Range("M85").Value = 445.5333000000001
$ws.Range("N85").Value = -4196
$ws.Range("H109").Value = 33349.5
$ws.Range("J109").Value = 33349.5
$ws.Range("L109").Value = 33349.5
$ws.Range("N109").Value = -36123.5
$ws.Range("H139").Value = 45426
$ws.Range("J139").Value = 45426
$ws.Range("L139").Value = 45426
$ws.Range("N139").Value = -55706

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3221.8276
$ws.Range("I136").Value = 511.5909
$ws.Range("J136").Value = 11739.714
$ws.Range("K136").Value = 1534.7727
$ws.Range("L136").Value = 35219.142
$ws.Range("M136").Value = 1015.2273
$ws.Range("N136").Value = -40319.142
$ws.Range("H139").Value = 39870.332
$ws.Range("J139").Value = 40195
$ws.Range("L139").Value = 40195
$ws.Range("N139").Value = -50475
